$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update F column values for several rows
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 24
$ws1.Range("F4").Value = 1380
$ws1.Range("F5").Value = 314
$ws1.Range("F7").Value = 10688
$ws1.Range("F8").Value = 22
$ws1.Range("F12").Value = 706
$ws1.Range("F13").Value = 12046
$ws1.Range("F14").Value = 12493

# Sheet "全部类型" (fourth sheet) - same data shifted down by one row
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 24
$ws4.Range("F5").Value = 1380
$ws4.Range("F6").Value = 314
$ws4.Range("F8").Value = 10688
$ws4.Range("F9").Value = 22
$ws4.Range("F13").Value = 706
$ws4.Range("F14").Value = 12046
$ws4.Range("F15").Value = 12493
